$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 20:52"

# Update country rows whose figures changed in this data refresh
# (values below represent the state after the 20:52 update; some countries
# swapped table position versus their neighbours as rankings shifted)
# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 418559
$ws.Range("C4").Value = 18224
$ws.Range("D4").Value = 22184
$ws.Range("E4").Value = 382135
$ws.Range("F4").Value = 9225
$ws.Range("G4").Value = 1399
$ws.Range("H4").Value = 14240

# Row 8: Alemania
$ws.Range("A8").Value = "Alemania"
$ws.Range("B8").Value = 110698
$ws.Range("C8").Value = 3035
$ws.Range("D8").Value = 36081
$ws.Range("E8").Value = 72425
$ws.Range("F8").Value = 4895
$ws.Range("G8").Value = 176
$ws.Range("H8").Value = 2192

# Row 17: Brasil
$ws.Range("A17").Value = "Brasil"
$ws.Range("B17").Value = 14511
$ws.Range("C17").Value = 477
$ws.Range("D17").Value = 127
$ws.Range("E17").Value = 13664
$ws.Range("F17").Value = 296
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 720

# Row 51: Sudafrica
$ws.Range("A51").Value = "Sudafrica"
$ws.Range("B51").Value = 1845
$ws.Range("C51").Value = 96
$ws.Range("D51").Value = 95
$ws.Range("E51").Value = 1732
$ws.Range("F51").Value = 7
$ws.Range("G51").Value = 5
$ws.Range("H51").Value = 18

# Row 52: Colombia
$ws.Range("A52").Value = "Colombia"
$ws.Range("B52").Value = 1780
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 100
$ws.Range("E52").Value = 1630
$ws.Range("F52").Value = 76
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 50

# Row 53: Argentina
$ws.Range("A53").Value = "Argentina"
$ws.Range("B53").Value = 1715
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 358
$ws.Range("E53").Value = 1294
$ws.Range("F53").Value = 96
$ws.Range("G53").Value = 3
$ws.Range("H53").Value = 63

# Row 92: Burkina Faso
$ws.Range("A92").Value = "Burkina Faso"
$ws.Range("B92").Value = 414
$ws.Range("C92").Value = 30
$ws.Range("D92").Value = 134
$ws.Range("E92").Value = 257
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 4
$ws.Range("H92").Value = 23

# Row 93: Albania
$ws.Range("A93").Value = "Albania"
$ws.Range("B93").Value = 400
$ws.Range("C93").Value = 17
$ws.Range("D93").Value = 154
$ws.Range("E93").Value = 224
$ws.Range("F93").Value = 7
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 22

# Row 135: Aruba
$ws.Range("A135").Value = "Aruba"
$ws.Range("B135").Value = 77
$ws.Range("C135").Value = 3
$ws.Range("D135").Value = 14
$ws.Range("E135").Value = 63
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 0

# Row 136: Guayana Francesa
$ws.Range("A136").Value = "Guayana Francesa"
$ws.Range("B136").Value = 77
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 34
$ws.Range("E136").Value = 43
$ws.Range("F136").Value = 1
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0

# Row 142: Uganda
$ws.Range("A142").Value = "Uganda"
$ws.Range("B142").Value = 53
$ws.Range("C142").Value = 1
$ws.Range("D142").Value = 0
$ws.Range("E142").Value = 53
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 0

# Row 143: Polinesia Francesa
$ws.Range("A143").Value = "Polinesia Francesa"
$ws.Range("B143").Value = 51
$ws.Range("C143").Value = 4
$ws.Range("D143").Value = 0
$ws.Range("E143").Value = 51
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 0

# Row 182: Suazilandia
$ws.Range("A182").Value = "Suazilandia"
$ws.Range("B182").Value = 12
$ws.Range("C182").Value = 2
$ws.Range("D182").Value = 7
$ws.Range("E182").Value = 5
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 0

# Row 185: Zimbabue
$ws.Range("A185").Value = "Zimbabue"
$ws.Range("B185").Value = 11
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 8
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 1
$ws.Range("H185").Value = 3

# Row 203: Islas Malvinas
$ws.Range("A203").Value = "Islas Malvinas"
$ws.Range("B203").Value = 5
$ws.Range("C203").Value = 3
$ws.Range("D203").Value = 1
$ws.Range("E203").Value = 4
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Row 204: Butan
$ws.Range("A204").Value = "Butan"
$ws.Range("B204").Value = 5
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 2
$ws.Range("E204").Value = 3
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

# Row 206: Sahara Occidental
$ws.Range("A206").Value = "Sahara Occidental"
$ws.Range("B206").Value = 4
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 0
$ws.Range("E206").Value = 4
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0

# Row 207: Gambia
$ws.Range("A207").Value = "Gambia"
$ws.Range("B207").Value = 4
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 2
$ws.Range("E207").Value = 1
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 1

# Row 210: Burundi
$ws.Range("A210").Value = "Burundi"
$ws.Range("B210").Value = 3
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 0
$ws.Range("E210").Value = 3
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 212: Papua Nueva Guinea
$ws.Range("A212").Value = "Papua Nueva Guinea"
$ws.Range("B212").Value = 2
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 0
$ws.Range("E212").Value = 2
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0

# Row 213: Bonaire, San Eustaquio y Saba
$ws.Range("A213").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B213").Value = 2
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 0
$ws.Range("E213").Value = 2
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214: San Pedro y Miquelon
$ws.Range("A214").Value = "San Pedro y Miquelon"
$ws.Range("B214").Value = 1
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 0
$ws.Range("E214").Value = 1
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# Row 215: Timor Oriental
$ws.Range("A215").Value = "Timor Oriental"
$ws.Range("B215").Value = 1
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 0
$ws.Range("E215").Value = 1
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

